$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2 values (identifier, title(blank), date_s(blank), levelOfDescription,
# extentAndMedium, notes, file_path(blank) -- per the sharedStrings/header layout)
$ws.Range("A2").Value = "MCH186-1"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 23O | GRAP COUNT NUMER: NONE"
$ws.Range("H2").Value = ""

# Match the font used for the new row of data: Calibri 10pt, automatic/theme text color
foreach ($addr in @("A2","C2","D2","E2","F2","G2","H2")) {
    $cell = $ws.Range($addr)
    $cell.Font.ThemeColor = 1
    $cell.Font.Name = "Calibri"
}
